# Auto-generated edit script: updates the crypto price table to the 28-1-2023 snapshot.
# Matches the commit "Updated symbol list on Sat Jan 28 00:53:18 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''308.72'
$ws.Range("E2").Value = '''1.37%'
$ws.Range("F2").Value = '28-1-2023'
$ws.Range("G2").Value = '''0'
$ws.Range("D3").Value = '''36.33'
$ws.Range("E3").Value = '''0.97%'
$ws.Range("F3").Value = '28-1-2023'
$ws.Range("G3").Value = '''0'
$ws.Range("E4").Value = '''0.53%'
$ws.Range("F4").Value = '28-1-2023'
$ws.Range("G4").Value = '''0'
$ws.Range("D5").Value = '''0.08186'
$ws.Range("E5").Value = '''2.08%'
$ws.Range("F5").Value = '28-1-2023'
$ws.Range("G5").Value = '''0'
$ws.Range("D6").Value = '''1.968'
$ws.Range("E6").Value = '''6.04%'
$ws.Range("F6").Value = '28-1-2023'
$ws.Range("G6").Value = '''0'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '''4.165'
$ws.Range("E7").Value = '''0.32%'
$ws.Range("F7").Value = '28-1-2023'
$ws.Range("G7").Value = '''0'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").Value = '''7.866'
$ws.Range("E8").Value = '''0.43%'
$ws.Range("F8").Value = '28-1-2023'
$ws.Range("G8").Value = '''0'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '''0.9319'
$ws.Range("E9").Value = '''0.37%'
$ws.Range("F9").Value = '28-1-2023'
$ws.Range("G9").Value = '''0'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '''0.1459'
$ws.Range("E10").Value = '''11.62%'
$ws.Range("F10").Value = '28-1-2023'
$ws.Range("G10").Value = '''0'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '''0.1941'
$ws.Range("E11").Value = '''1.82%'
$ws.Range("F11").Value = '28-1-2023'
$ws.Range("G11").Value = '''0'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '''0.09145'
$ws.Range("E12").Value = '''-1.15%'
$ws.Range("F12").Value = '28-1-2023'
$ws.Range("G12").Value = '''0'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '''0.03510'
$ws.Range("E13").Value = '''-0.31%'
$ws.Range("F13").Value = '28-1-2023'
$ws.Range("G13").Value = '''0'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '''0.09870'
$ws.Range("E14").Value = '''-0.13%'
$ws.Range("F14").Value = '28-1-2023'
$ws.Range("G14").Value = '''0'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '''0.001400'
$ws.Range("E15").Value = '''-1.25%'
$ws.Range("F15").Value = '28-1-2023'
$ws.Range("G15").Value = '''0'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '''0.006681'
$ws.Range("E16").Value = '''6.00%'
$ws.Range("F16").Value = '28-1-2023'
$ws.Range("G16").Value = '''0'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '''3.859'
$ws.Range("E17").Value = '''5.50%'
$ws.Range("F17").Value = '28-1-2023'
$ws.Range("G17").Value = '''0'
$ws.Range("D18").Value = '''3.438'
$ws.Range("E18").Value = '''4.90%'
$ws.Range("F18").Value = '28-1-2023'
$ws.Range("G18").Value = '''0'
$ws.Range("D19").Value = '''0.3453'
$ws.Range("E19").Value = '''0.06%'
$ws.Range("F19").Value = '28-1-2023'
$ws.Range("G19").Value = '''0'
$ws.Range("D20").Value = '''0.1333'
$ws.Range("E20").Value = '''-1.21%'
$ws.Range("F20").Value = '28-1-2023'
$ws.Range("G20").Value = '''0'
$ws.Range("D21").Value = '''4.817'
$ws.Range("E21").Value = '''-7.33%'
$ws.Range("F21").Value = '28-1-2023'
$ws.Range("G21").Value = '''0'
$ws.Range("D22").Value = '''0.2611'
$ws.Range("E22").Value = '''2.91%'
$ws.Range("F22").Value = '28-1-2023'
$ws.Range("G22").Value = '''0'
$ws.Range("D23").Value = '''0.04386'
$ws.Range("E23").Value = '''-1.26%'
$ws.Range("F23").Value = '28-1-2023'
$ws.Range("G23").Value = '''0'
$ws.Range("D24").Value = '''0.001233'
$ws.Range("E24").Value = '''-0.34%'
$ws.Range("F24").Value = '28-1-2023'
$ws.Range("G24").Value = '''0'
$ws.Range("D25").Value = '''0.004174'
$ws.Range("E25").Value = '''-11.09%'
$ws.Range("F25").Value = '28-1-2023'
$ws.Range("G25").Value = '''0'
$ws.Range("F26").Value = '28-1-2023'
$ws.Range("G26").Value = '''0'
$ws.Range("D27").Value = '''0.0001300'
$ws.Range("E27").Value = '''-0.21%'
$ws.Range("F27").Value = '28-1-2023'
$ws.Range("G27").Value = '''0'
$ws.Range("F28").Value = '28-1-2023'
$ws.Range("G28").Value = '''0'
$ws.Range("F29").Value = '28-1-2023'
$ws.Range("G29").Value = '''0'
$ws.Range("F30").Value = '28-1-2023'
$ws.Range("G30").Value = '''0'
$ws.Range("F31").Value = '28-1-2023'
$ws.Range("G31").Value = '''0'
$ws.Range("F32").Value = '28-1-2023'
$ws.Range("G32").Value = '''0'
$ws.Range("F33").Value = '28-1-2023'
$ws.Range("G33").Value = '''0'
$ws.Range("F34").Value = '28-1-2023'
$ws.Range("G34").Value = '''0'
$ws.Range("F35").Value = '28-1-2023'
$ws.Range("G35").Value = '''0'
$ws.Range("F36").Value = '28-1-2023'
$ws.Range("G36").Value = '''0'
$ws.Range("F37").Value = '28-1-2023'
$ws.Range("G37").Value = '''0'
$ws.Range("F38").Value = '28-1-2023'
$ws.Range("G38").Value = '''0'
$ws.Range("D39").Value = '''0.02098'
$ws.Range("E39").Value = '''8.09%'
$ws.Range("F39").Value = '28-1-2023'
$ws.Range("G39").Value = '''0'
$ws.Range("D40").Value = '''0.05158'
$ws.Range("E40").Value = '''0.28%'
$ws.Range("F40").Value = '28-1-2023'
$ws.Range("G40").Value = '''0'
$ws.Range("D41").Value = '''0.007477'
$ws.Range("E41").Value = '''-1.14%'
$ws.Range("F41").Value = '28-1-2023'
$ws.Range("G41").Value = '''0'
$ws.Range("D42").Value = '''0.01007'
$ws.Range("E42").Value = '''-1.26%'
$ws.Range("F42").Value = '28-1-2023'
$ws.Range("G42").Value = '''0'
$ws.Range("D43").Value = '''0.1372'
$ws.Range("E43").Value = '''0.40%'
$ws.Range("F43").Value = '28-1-2023'
$ws.Range("G43").Value = '''0'
$ws.Range("D44").Value = '''0.002131'
$ws.Range("E44").Value = '''-1.14%'
$ws.Range("F44").Value = '28-1-2023'
$ws.Range("G44").Value = '''0'
$ws.Range("D45").Value = '''0.01004'
$ws.Range("E45").Value = '''1.70%'
$ws.Range("F45").Value = '28-1-2023'
$ws.Range("G45").Value = '''0'
$ws.Range("D46").Value = '''0.00006380'
$ws.Range("E46").Value = '''0.99%'
$ws.Range("F46").Value = '28-1-2023'
$ws.Range("G46").Value = '''0'
$ws.Range("D47").Value = '''0.00000000750'
$ws.Range("E47").Value = '''-0.21%'
$ws.Range("F47").Value = '28-1-2023'
$ws.Range("G47").Value = '''0'
$ws.Range("E48").Value = '''-0.24%'
$ws.Range("F48").Value = '28-1-2023'
$ws.Range("G48").Value = '''0'
$ws.Range("D49").Value = '''0.001598'
$ws.Range("E49").Value = '''-3.93%'
$ws.Range("F49").Value = '28-1-2023'
$ws.Range("G49").Value = '''0'
$ws.Range("D50").Value = '''0.00002100'
$ws.Range("E50").Value = '''-0.21%'
$ws.Range("F50").Value = '28-1-2023'
$ws.Range("G50").Value = '''0'
$ws.Range("D51").Value = '''0.0002000'
$ws.Range("E51").Value = '''-0.21%'
$ws.Range("F51").Value = '28-1-2023'
$ws.Range("G51").Value = '''0'
